$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data row (row 10) -------------------------------------------------
$ws.Range("B10").Value = 20
$ws.Range("C10").Formula = "=3*29054"
$ws.Range("D10").Formula = "=C10*F10"
$ws.Range("E10").Value = 389
$ws.Range("F10").Value = 2
$ws.Range("G10").Value = 1366
$ws.Range("H10").NumberFormat = "h:mm:ss"
$ws.Range("H10").Value = 0.03439814814814815
$ws.Range("I10").Value = 6800
$ws.Range("J10").Value = "Vampiro"
$ws.Range("K10").Value = "Normal"
$ws.Range("L10").NumberFormat = "m/d/yy"
$ws.Range("L10").Value = 46014

# Update the active selection to match the new edit location
$ws.Range("N9").Select() | Out-Null
